{"js": "/*\n * Applies the table-cell value replacements described by the diff.\n * The document has a single 20x5 table (100 cells), each holding one\n * arithmetic equation as plain text (e.g. \"76-31=45\"). The diff replaces\n * each cell's equation text with a new one, in row-major document order.\n * We address cells positionally (table.getCell(row, col)) rather than by\n * searching for the old text, because a couple of old values repeat\n * (e.g. \"8+68=76\" occurs twice but maps to two different replacements).\n */\n\n// Old (pre-edit) and new (post-edit) cell values, in row-major order\n// (row 0 col 0, row 0 col 1, ..., row 0 col 4, row 1 col 0, ...).\nconst oldValues = [\n  \"76-31=45\",\n  \"52-43=9\",\n  \"49-39=10\",\n  \"32+10=42\",\n  \"32-5=27\",\n  \"74-45=29\",\n  \"33-24=9\",\n  \"69-43=26\",\n  \"45-4=41\",\n  \"57-47=10\",\n  \"84-44=40\",\n  \"67-2=65\",\n  \"2+4=6\",\n  \"79-16=63\",\n  \"21-8=13\",\n  \"8+68=76\",\n  \"17+0=17\",\n  \"58+13=71\",\n  \"38+9=47\",\n  \"85-43=42\",\n  \"13-9=4\",\n  \"14+58=72\",\n  \"78-10=68\",\n  \"83-68=15\",\n  \"47+9=56\",\n  \"26+50=76\",\n  \"6+37=43\",\n  \"86+13=99\",\n  \"11+35=46\",\n  \"23+70=93\",\n  \"64+23=87\",\n  \"17+82=99\",\n  \"41+42=83\",\n  \"88-27=61\",\n  \"52+35=87\",\n  \"89-15=74\",\n  \"4+94=98\",\n  \"9+19=28\",\n  \"81-71=10\",\n  \"94-57=37\",\n  \"91-20=71\",\n  \"98-35=63\",\n  \"88-18=70\",\n  \"46-44=2\",\n  \"63-24=39\",\n  \"51+19=70\",\n  \"41+58=99\",\n  \"48-41=7\",\n  \"25-10=15\",\n  \"33+56=89\",\n  \"11+5=16\",\n  \"96-62=34\",\n  \"13+73=86\",\n  \"71-59=12\",\n  \"97-77=20\",\n  \"70-63=7\",\n  \"82-82=0\",\n  \"22-18=4\",\n  \"11+36=47\",\n  \"48+49=97\",\n  \"82-18=64\",\n  \"25-16=9\",\n  \"23+60=83\",\n  \"38+54=92\",\n  \"8+68=76\",\n  \"24+1=25\",\n  \"83-25=58\",\n  \"58-8=50\",\n  \"35-15=20\",\n  \"80+10=90\",\n  \"45-30=15\",\n  \"24+59=83\",\n  \"46-13=33\",\n  \"57-25=32\",\n  \"43+22=65\",\n  \"54+37=91\",\n  \"32+55=87\",\n  \"50-0=50\",\n  \"15+73=88\",\n  \"83-43=40\",\n  \"8+65=73\",\n  \"55-5=50\",\n  \"51+12=63\",\n  \"67+6=73\",\n  \"52-30=22\",\n  \"86-30=56\",\n  \"66+30=96\",\n  \"9+55=64\",\n  \"51+33=84\",\n  \"72-28=44\",\n  \"87-85=2\",\n  \"52+28=80\",\n  \"74-47=27\",\n  \"94-91=3\",\n  \"32+23=55\",\n  \"63-9=54\",\n  \"47+47=94\",\n  \"66+25=91\",\n  \"33+53=86\",\n  \"63+34=97\"\n];\nconst newValues = [\n  \"16+79=95\",\n  \"52+39=91\",\n  \"45-34=11\",\n  \"10+70=80\",\n  \"99-40=59\",\n  \"69+6=75\",\n  \"18+16=34\",\n  \"66-2=64\",\n  \"51+47=98\",\n  \"47+49=96\",\n  \"59-19=40\",\n  \"85-60=25\",\n  \"94-17=77\",\n  \"84-80=4\",\n  \"74-60=14\",\n  \"74-33=41\",\n  \"50+21=71\",\n  \"2+77=79\",\n  \"42+53=95\",\n  \"20-13=7\",\n  \"39-33=6\",\n  \"36-2=34\",\n  \"77-47=30\",\n  \"12+30=42\",\n  \"52-26=26\",\n  \"22-21=1\",\n  \"90-9=81\",\n  \"12+15=27\",\n  \"88-61=27\",\n  \"9+80=89\",\n  \"27+58=85\",\n  \"51+4=55\",\n  \"18+40=58\",\n  \"21+55=76\",\n  \"69+22=91\",\n  \"82-52=30\",\n  \"32+13=45\",\n  \"31-4=27\",\n  \"54+33=87\",\n  \"16+5=21\",\n  \"3+61=64\",\n  \"19+8=27\",\n  \"90-48=42\",\n  \"21+62=83\",\n  \"63-41=22\",\n  \"89+6=95\",\n  \"61-43=18\",\n  \"20+78=98\",\n  \"56+0=56\",\n  \"42+5=47\",\n  \"26-2=24\",\n  \"40+35=75\",\n  \"84-8=76\",\n  \"87-40=47\",\n  \"74+23=97\",\n  \"44-31=13\",\n  \"2+8=10\",\n  \"89-40=49\",\n  \"46+34=80\",\n  \"84-13=71\",\n  \"92-37=55\",\n  \"11+50=61\",\n  \"99-38=61\",\n  \"26+58=84\",\n  \"68-64=4\",\n  \"25-6=19\",\n  \"87-80=7\",\n  \"35-28=7\",\n  \"18+44=62\",\n  \"99-84=15\",\n  \"39+12=51\",\n  \"60-10=50\",\n  \"23-10=13\",\n  \"21+23=44\",\n  \"23+54=77\",\n  \"1+76=77\",\n  \"70-44=26\",\n  \"12+13=25\",\n  \"41-28=13\",\n  \"77-34=43\",\n  \"22-8=14\",\n  \"58+2=60\",\n  \"54+10=64\",\n  \"88-37=51\",\n  \"76-73=3\",\n  \"5+28=33\",\n  \"78+20=98\",\n  \"72-72=0\",\n  \"9+70=79\",\n  \"71-65=6\",\n  \"97-40=57\",\n  \"11+64=75\",\n  \"55+20=75\",\n  \"57-44=13\",\n  \"52+13=65\",\n  \"67-63=4\",\n  \"21+8=29\",\n  \"49-24=25\",\n  \"37+38=75\",\n  \"27+44=71\"\n];\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst columnCount = 5;\nconst totalCells = oldValues.length;\n\n// Sanity-check: load every current cell value first, so we can confirm we\n// are editing the expected cells before mutating anything.\nconst cells = [];\nfor (let i = 0; i < totalCells; i++) {\n  const row = Math.floor(i / columnCount);\n  const col = i % columnCount;\n  const cell = table.getCell(row, col);\n  cell.load(\"value\");\n  cells.push(cell);\n}\nawait context.sync();\n\nfor (let i = 0; i < totalCells; i++) {\n  const expected = oldValues[i];\n  const actual = cells[i].value;\n  if (actual !== expected) {\n    throw new Error(\n      \"Unexpected cell value at index \" + i + \": expected '\" + expected +\n      \"' but found '\" + actual + \"'\"\n    );\n  }\n}\n\n// Now apply the replacements.\nfor (let i = 0; i < totalCells; i++) {\n  cells[i].value = newValues[i];\n}\n\nawait context.sync();\n", "ps1": "# Applies the table-cell value replacements described by the diff.\n# The document has a single 20x5 table (100 cells), each holding one\n# arithmetic equation as plain text (e.g. \"76-31=45\"). The diff replaces\n# each cell's equation text with a new one, in row-major document order.\n# We address cells positionally (Table.Cell(row, col), 1-based) rather than\n# by searching for the old text, because a couple of old values repeat\n# (e.g. \"8+68=76\" occurs twice but maps to two different replacements).\n\n$d = $word.ActiveDocument\n\n# Old (pre-edit) and new (post-edit) cell values, in row-major order\n# (row 1 col 1, row 1 col 2, ..., row 1 col 5, row 2 col 1, ...).\n$oldValues = @(\n    '76-31=45',\n    '52-43=9',\n    '49-39=10',\n    '32+10=42',\n    '32-5=27',\n    '74-45=29',\n    '33-24=9',\n    '69-43=26',\n    '45-4=41',\n    '57-47=10',\n    '84-44=40',\n    '67-2=65',\n    '2+4=6',\n    '79-16=63',\n    '21-8=13',\n    '8+68=76',\n    '17+0=17',\n    '58+13=71',\n    '38+9=47',\n    '85-43=42',\n    '13-9=4',\n    '14+58=72',\n    '78-10=68',\n    '83-68=15',\n    '47+9=56',\n    '26+50=76',\n    '6+37=43',\n    '86+13=99',\n    '11+35=46',\n    '23+70=93',\n    '64+23=87',\n    '17+82=99',\n    '41+42=83',\n    '88-27=61',\n    '52+35=87',\n    '89-15=74',\n    '4+94=98',\n    '9+19=28',\n    '81-71=10',\n    '94-57=37',\n    '91-20=71',\n    '98-35=63',\n    '88-18=70',\n    '46-44=2',\n    '63-24=39',\n    '51+19=70',\n    '41+58=99',\n    '48-41=7',\n    '25-10=15',\n    '33+56=89',\n    '11+5=16',\n    '96-62=34',\n    '13+73=86',\n    '71-59=12',\n    '97-77=20',\n    '70-63=7',\n    '82-82=0',\n    '22-18=4',\n    '11+36=47',\n    '48+49=97',\n    '82-18=64',\n    '25-16=9',\n    '23+60=83',\n    '38+54=92',\n    '8+68=76',\n    '24+1=25',\n    '83-25=58',\n    '58-8=50',\n    '35-15=20',\n    '80+10=90',\n    '45-30=15',\n    '24+59=83',\n    '46-13=33',\n    '57-25=32',\n    '43+22=65',\n    '54+37=91',\n    '32+55=87',\n    '50-0=50',\n    '15+73=88',\n    '83-43=40',\n    '8+65=73',\n    '55-5=50',\n    '51+12=63',\n    '67+6=73',\n    '52-30=22',\n    '86-30=56',\n    '66+30=96',\n    '9+55=64',\n    '51+33=84',\n    '72-28=44',\n    '87-85=2',\n    '52+28=80',\n    '74-47=27',\n    '94-91=3',\n    '32+23=55',\n    '63-9=54',\n    '47+47=94',\n    '66+25=91',\n    '33+53=86',\n    '63+34=97'\n)\n$newValues = @(\n    '16+79=95',\n    '52+39=91',\n    '45-34=11',\n    '10+70=80',\n    '99-40=59',\n    '69+6=75',\n    '18+16=34',\n    '66-2=64',\n    '51+47=98',\n    '47+49=96',\n    '59-19=40',\n    '85-60=25',\n    '94-17=77',\n    '84-80=4',\n    '74-60=14',\n    '74-33=41',\n    '50+21=71',\n    '2+77=79',\n    '42+53=95',\n    '20-13=7',\n    '39-33=6',\n    '36-2=34',\n    '77-47=30',\n    '12+30=42',\n    '52-26=26',\n    '22-21=1',\n    '90-9=81',\n    '12+15=27',\n    '88-61=27',\n    '9+80=89',\n    '27+58=85',\n    '51+4=55',\n    '18+40=58',\n    '21+55=76',\n    '69+22=91',\n    '82-52=30',\n    '32+13=45',\n    '31-4=27',\n    '54+33=87',\n    '16+5=21',\n    '3+61=64',\n    '19+8=27',\n    '90-48=42',\n    '21+62=83',\n    '63-41=22',\n    '89+6=95',\n    '61-43=18',\n    '20+78=98',\n    '56+0=56',\n    '42+5=47',\n    '26-2=24',\n    '40+35=75',\n    '84-8=76',\n    '87-40=47',\n    '74+23=97',\n    '44-31=13',\n    '2+8=10',\n    '89-40=49',\n    '46+34=80',\n    '84-13=71',\n    '92-37=55',\n    '11+50=61',\n    '99-38=61',\n    '26+58=84',\n    '68-64=4',\n    '25-6=19',\n    '87-80=7',\n    '35-28=7',\n    '18+44=62',\n    '99-84=15',\n    '39+12=51',\n    '60-10=50',\n    '23-10=13',\n    '21+23=44',\n    '23+54=77',\n    '1+76=77',\n    '70-44=26',\n    '12+13=25',\n    '41-28=13',\n    '77-34=43',\n    '22-8=14',\n    '58+2=60',\n    '54+10=64',\n    '88-37=51',\n    '76-73=3',\n    '5+28=33',\n    '78+20=98',\n    '72-72=0',\n    '9+70=79',\n    '71-65=6',\n    '97-40=57',\n    '11+64=75',\n    '55+20=75',\n    '57-44=13',\n    '52+13=65',\n    '67-63=4',\n    '21+8=29',\n    '49-24=25',\n    '37+38=75',\n    '27+44=71'\n)\n\n$table = $d.Tables.Item(1)\n$columnCount = 5\n$totalCells = $oldValues.Count\n\n# Sanity-check every current cell value before mutating anything, so we\n# only touch the cells we expect to.\nfor ($i = 0; $i -lt $totalCells; $i++) {\n    $row = [int][Math]::Floor($i / $columnCount) + 1\n    $col = ($i % $columnCount) + 1\n    $cell = $table.Cell($row, $col)\n    # Cell.Range.Text includes the trailing cell-mark (\\r\\a); strip it off\n    # before comparing against the plain equation text.\n    $actual = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    $expected = $oldValues[$i]\n    if ($actual -ne $expected) {\n        throw \"Unexpected cell value at index $i (row $row, col $col): expected '$expected' but found '$actual'\"\n    }\n}\n\n# Now apply the replacements.\nfor ($i = 0; $i -lt $totalCells; $i++) {\n    $row = [int][Math]::Floor($i / $columnCount) + 1\n    $col = ($i % $columnCount) + 1\n    $cell = $table.Cell($row, $col)\n    $cell.Range.Text = $newValues[$i]\n}\n"}
